$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.947.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").Value = "'1.876.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.96%  "

$ws.Range("D5").Value = "'313.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D7").Value = "'0.4861"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "

$ws.Range("D8").Value = "'0.3813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "

$ws.Range("D9").Value = "'0.07370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Value = "'0.9417"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("D11").Value = "'21.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.28%  "

$ws.Range("D12").Value = "'0.07778"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").Value = "'1.908.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("E14").Value = "  +2.33%  "

$ws.Range("D15").Value = "'6.602"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").Value = "'91.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").Value = "'1.012"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").Value = "'0.000008878"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("D20").Value = "'27.955.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "

$ws.Range("D21").Value = "'14.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "'5.123"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").Value = "'2.141.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "

$ws.Range("D24").Value = "'10.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.949"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'157.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.45%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "'2.043"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.25%  "

$ws.Range("D29").Value = "'115.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").Value = "'4.972"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("D31").Value = "'0.08902"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = "'3.338"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").Value = "'1.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").Value = "'0.7716"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.15%  "

$ws.Range("D35").Value = "'4.645"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").Value = "'2.727"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.57%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").Value = "'0.02047"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'0.5599"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.16%  "

$ws.Range("D40").Value = "'0.05375"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("D41").Value = "'3.006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").Value = "'8.538"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").Value = "'0.1520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").Value = "'0.4890"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("D46").Value = "'10.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("D47").Value = "'105.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").Value = "'1.670"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.16%  "

$ws.Range("D50").Value = "'68.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "'0.06108"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "

Write-Host "Applied cryptos update"
